$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.975.30'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.119.66'
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.36'
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.26'
$ws.Range("E6").Value = '  +2.06%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.39'
$ws.Range("E9").Value = '  -3.21%  '
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.15'
$ws.Range("E13").Value = '  +1.86%  '
$ws.Range("E14").Value = '  -1.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.636.15'
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.923.46'
$ws.Range("E16").Value = '  +0.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.15'
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.119.72'
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.26'
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '475.24'
$ws.Range("E20").Value = '  +1.49%  '
$ws.Range("E21").Value = '  -0.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.88'
$ws.Range("E22").Value = '  +5.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '83.84'
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("E25").Value = '  -3.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.27'
$ws.Range("E26").Value = '  +1.48%  '
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.91'
$ws.Range("E28").Value = '  -1.72%  '
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.63'
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0951'
$ws.Range("E33").Value = '  -6.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.84'
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("E36").Value = '  -2.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '46.79'
$ws.Range("E37").Value = '  -1.16%  '
$ws.Range("E38").Value = '  -0.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.06'
$ws.Range("E39").Value = '  -1.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.312'
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("E41").Value = '  +1.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.61'
$ws.Range("E42").Value = '  -1.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.819.55'
$ws.Range("E43").Value = '  +1.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '383.12'
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0353'
$ws.Range("E45").Value = '  -2.38%  '
$ws.Range("E46").Value = '  -9.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '135.89'
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.96'
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("E51").Value = '  -0.81%  '
